# Auto-generated edit script applying cryptos.xlsx price/volume updates
# (Sat Nov 11 21:51:54 UTC 2023 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.024.18"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.053.01"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.02"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.669"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.37"
$ws.Range("E7").Value = "  +8.09%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.15"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").Value = "  +4.99%  "
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.35"
$ws.Range("E13").Value = "  +7.58%  "
$ws.Range("D14").Value = "2.352.35"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.807"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.60"
$ws.Range("E16").Value = "  +6.91%  "
$ws.Range("D17").Value = "2.063.09"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "37.005.87"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.82"
$ws.Range("E19").Value = "  +15.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.25"
$ws.Range("E20").Value = "  +2.78%  "
$ws.Range("E21").Value = "  +6.49%  "
$ws.Range("E22").Value = "  +3.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.88"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  -3.71%  "
$ws.Range("E26").Value = "  +11.72%  "
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.30"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.24"
$ws.Range("E29").Value = "  -3.40%  "
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  +5.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.75"
$ws.Range("E32").Value = "  +4.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.47"
$ws.Range("E34").Value = "  +5.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0890"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("E38").Value = "  -3.23%  "
$ws.Range("E39").Value = "  +19.80%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.73"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0225"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.43"
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.60"
$ws.Range("E46").Value = "  +12.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("E47").Value = "  +4.55%  "
$ws.Range("D48").Value = "1.289.71"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.67"
$ws.Range("E51").Value = "  -25.95%  "
